$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ECs" target-cluster rows (original rows 8,9,10, where Target cluster = ECs)
# are removed entirely, and the remaining rows' derived-specificity statistics are
# recomputed now that "ECs" is no longer one of the possible target clusters.
$ws.Rows("8:10").Delete()

# Row 2: Sending=ECs, Target=FAPs (was ECs)
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 0.5347833333333334
$ws.Range("H2").Value = 1.60435
$ws.Range("I2").Value = 0.196822066153855
$ws.Range("J2").Value = 0.196822066153855
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.183046666666667
$ws.Range("N2").Value = 3.54914
$ws.Range("O2").Value = 0.6222589862820888
$ws.Range("P2").Value = 0.6222589862820888
$ws.Range("Q2").Value = 0.6326736398888889
$ws.Range("R2").Value = 5.694062759
$ws.Range("S2").Value = 0.122474299362844
$ws.Range("T2").Value = 0.122474299362844

# Row 3: Sending=ECs, Target=MuSCs (was FAPs)
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 0.5347833333333334
$ws.Range("H3").Value = 1.60435
$ws.Range("I3").Value = 0.196822066153855
$ws.Range("J3").Value = 0.196822066153855
$ws.Range("M3").Value = 0.718166
$ws.Range("N3").Value = 2.154498
$ws.Range("O3").Value = 0.3777410137179113
$ws.Range("P3").Value = 0.3777410137179112
$ws.Range("Q3").Value = 0.3840632073666667
$ws.Range("R3").Value = 3.4565688663
$ws.Range("S3").Value = 0.074347766791011
$ws.Range("T3").Value = 0.07434776679101097

# Row 4: Sending=FAPs (was ECs), Target=FAPs (was MuSCs)
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.514066
$ws.Range("H4").Value = 1.542198
$ws.Range("I4").Value = 0.1891972429821067
$ws.Range("J4").Value = 0.1891972429821067
$ws.Range("M4").Value = 1.183046666666667
$ws.Range("N4").Value = 3.54914
$ws.Range("O4").Value = 0.6222589862820888
$ws.Range("P4").Value = 0.6222589862820888
$ws.Range("Q4").Value = 0.6081640677466666
$ws.Range("R4").Value = 5.47347660972
$ws.Range("S4").Value = 0.1177296846254118
$ws.Range("T4").Value = 0.1177296846254117

# Row 5: Sending=FAPs, Target=MuSCs (was ECs)
$ws.Range("D5").Value = "MuSCs"
$ws.Range("I5").Value = 0.1891972429821067
$ws.Range("J5").Value = 0.1891972429821067
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.718166
$ws.Range("N5").Value = 2.154498
$ws.Range("O5").Value = 0.3777410137179113
$ws.Range("P5").Value = 0.3777410137179112
$ws.Range("Q5").Value = 0.369184722956
$ws.Range("R5").Value = 3.322662506603999
$ws.Range("S5").Value = 0.07146755835669497
$ws.Range("T5").Value = 0.07146755835669495

# Row 6: Sending=MuSCs (was FAPs), Target=FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.668241
$ws.Range("H6").Value = 5.004723
$ws.Range("I6").Value = 0.6139806908640383
$ws.Range("J6").Value = 0.6139806908640382
$ws.Range("O6").Value = 0.6222589862820888
$ws.Range("P6").Value = 0.6222589862820888
$ws.Range("Q6").Value = 1.973606954246667
$ws.Range("R6").Value = 17.76246258822
$ws.Range("S6").Value = 0.382055002293833
$ws.Range("T6").Value = 0.3820550022938329

# Row 7: Sending=MuSCs (was FAPs), Target=MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.668241
$ws.Range("H7").Value = 5.004723
$ws.Range("I7").Value = 0.6139806908640383
$ws.Range("J7").Value = 0.6139806908640382
$ws.Range("O7").Value = 0.3777410137179113
$ws.Range("P7").Value = 0.3777410137179112
$ws.Range("Q7").Value = 1.198073966006
$ws.Range("R7").Value = 10.782665694054
$ws.Range("S7").Value = 0.2319256885702053
$ws.Range("T7").Value = 0.2319256885702053
